# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet (cloned from "2022-Q3" so it inherits the
#    exact same styling/layout) positioned right before "2022-Q3", and fill
#    it in with the new quarter's single fund row.
# 2. Insert a new row at the top of the "总计" (totals) sheet's data table
#    for the 2022-Q4 summary figures, pushing the existing 2022-Q3 / 2022-Q2
#    rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" worksheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The source sheet had two fund rows; 2022-Q4 only has one, so drop row 3.
$q4.Rows.Item(3).Delete()

# Store the numeric-looking figures as text (matching the source data),
# leaving A2 (index) and H2 (rank) as real numbers.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Cells.Item(2, 2).Value = "519029"
$q4.Cells.Item(2, 3).Value = "华夏稳增混合"
$q4.Cells.Item(2, 4).Value = "9.01"
$q4.Cells.Item(2, 5).Value = "94.55"
$q4.Cells.Item(2, 6).Value = "5.97"
$q4.Cells.Item(2, 7).Value = "0.5379"
$q4.Cells.Item(2, 8).Value = 3

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q4 summary row to the "总计" sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# Re-apply the index-column style (bold/bordered) that row 2 should carry,
# by copying it from the row directly below (now the old "2022-Q3" row).
$totals.Cells.Item(3, 1).Copy()
$totals.Cells.Item(2, 1).PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q4"
$totals.Cells.Item(2, 3).Value = 1
$totals.Cells.Item(2, 4).Value = 0.54

# The index column (A) is a 0-based running count; renumber the rows that
# shifted down so it stays sequential (0, 1, 2, ...).
$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(4, 1).Value = 2

# Restore the originally-active sheet ("2022-Q2") as the selected tab —
# cloning "2022-Q3" above leaves the new "2022-Q4" sheet active instead.
$wb.Worksheets.Item("2022-Q2").Activate()
